$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - worksheet index 1
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value = 7315
$wsExhibit.Range("F7").Value = 203
$wsExhibit.Range("F8").Value = 142
$wsExhibit.Range("F11").Value = 65
$wsExhibit.Range("F12").Value = 232
$wsExhibit.Range("F13").Value = 20
$wsExhibit.Range("F14").Value = 474
$wsExhibit.Range("F16").Value = 1878
$wsExhibit.Range("F17").Value = 52
$wsExhibit.Range("F19").Value = 3842
$wsExhibit.Range("F21").Value = 254
$wsExhibit.Range("F23").Value = 52
$wsExhibit.Range("F25").Value = 42
$wsExhibit.Range("F26").Value = 2528
$wsExhibit.Range("F28").Value = 333
$wsExhibit.Range("F30").Value = 9
$wsExhibit.Range("F33").Value = 31
$wsExhibit.Range("F37").Value = 170
$wsExhibit.Range("F38").Value = 54
$wsExhibit.Range("F39").Value = 1513
$wsExhibit.Range("F40").Value = 181

# Sheet "全部类型" (sheet4, all-types aggregate) - worksheet index 4
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 7315
$wsAll.Range("F8").Value = 203
$wsAll.Range("F9").Value = 142
$wsAll.Range("F12").Value = 65
$wsAll.Range("F13").Value = 232
$wsAll.Range("F14").Value = 20
$wsAll.Range("F15").Value = 474
$wsAll.Range("F17").Value = 1878
$wsAll.Range("F18").Value = 52
$wsAll.Range("F20").Value = 3842
$wsAll.Range("F22").Value = 254
$wsAll.Range("F24").Value = 52
$wsAll.Range("F26").Value = 42
$wsAll.Range("F27").Value = 2528
$wsAll.Range("F29").Value = 333
$wsAll.Range("F31").Value = 9
$wsAll.Range("F34").Value = 31
$wsAll.Range("F38").Value = 170
$wsAll.Range("F39").Value = 54
$wsAll.Range("F40").Value = 1513
$wsAll.Range("F41").Value = 181

